$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

$ws.Range("A3").Value = "ISTAT_29_7"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "ISTAT_31_739"
$ws.Range("B4").Value = 0

$ws.Range("A5").Select()
